$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 11: new vocabulary entry "intermission" ----

# A11: headword
$ws.Range("A11").Value = "brief intermission"

# B11: IPA phonetic transcription, formatted as alternating runs of
# Tahoma (for the IPA glyphs) and the default 等线 font (for the plain
# latin letters), matching the style used by the other phonetic cells
# in the sheet.
$ph = $ws.Range("B11")
$ph.Value = "[ˌɪntərˈmɪʃn]"
$ph.VerticalAlignment = -4108

$ph.Characters(2, 2).Font.Name = "Tahoma"
$ph.Characters(4, 2).Font.Name = "等线"
$ph.Characters(6, 1).Font.Name = "Tahoma"
$ph.Characters(7, 1).Font.Name = "等线"
$ph.Characters(8, 1).Font.Name = "Tahoma"
$ph.Characters(9, 1).Font.Name = "等线"
$ph.Characters(10, 2).Font.Name = "Tahoma"
$ph.Characters(12, 2).Font.Name = "等线"

# C11: definition
$ws.Range("C11").Value = "N-COUNT An intermission is a short break between two parts of a concert, show, or film. 幕间休息；中场休息"
$ws.Range("C11").WrapText = $true

# Cursor position left where the author's Excel session ended up.
[void]$ws.Range("C15").Select()

"done"
